# Adds the required "experimental" boolean element to the ValueSet metadata
# sheet (B7, next to "Experimental" in A7) and updates the Date value (B8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 must hold the literal text "true" (not an Excel TRUE boolean value).
# Assigning the bare string gets auto-coerced to a boolean by the COM
# layer, so stage it in a scratch cell with a leading apostrophe (forces
# text) and paste the value over - this keeps the normal "s=2" data style
# instead of picking up a stray quote-prefix style.
$ws.Range("Z1").Value = "'true"
$ws.Range("Z1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# The "Date" property value is refreshed to the new generation timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
